$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Total" row (row 42) had a blank leading cell (A42) followed by its
# data in B42:E42. Shift that row's contents one column to the left -
# equivalent to selecting A42 and doing Home > Delete > Delete Cells >
# Shift Cells Left - so the data now starts in A42 and the trailing cell
# (E42) becomes blank.
$ws.Range("B42:E42").Copy($ws.Range("A42"))
$ws.Range("E42").Value = ""
$ws.Range("E42").Style = "Normal"

# Update the active cell / selection to A42.
$ws.Range("A42").Select()
